# Update gh-pages output figures (想去人数 / "want to go" counts) across
# the four sheets, and refresh one event's title + cover image URL.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (Exhibitions) ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 455
$ws.Range("F5").Value = 8757
$ws.Range("F6").Value = 16
$ws.Range("F7").Value = 11165
$ws.Range("F18").Value = 86
$ws.Range("F22").Value = 1903
$ws.Range("F23").Value = 710
$ws.Range("F24").Value = 635
$ws.Range("F25").Value = 359
$ws.Range("F30").Value = 1303
$ws.Range("F39").Value = 307
$ws.Range("F43").Value = 389
$ws.Range("F44").Value = 118
$ws.Range("C45").Value = "北京·代号鸢only同人展（如鸢同人only）"
$ws.Range("F45").Value = 814
$ws.Range("I45").Value = "//i1.hdslb.com/bfs/openplatform/202409/p7GKpOhb1727406520811.jpeg"
$ws.Range("F48").Value = 158
$ws.Range("F49").Value = 144

# ---- Sheet "演出" (Performances) ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F15").Value = 13
$ws.Range("F16").Value = 14
$ws.Range("F24").Value = 395

# ---- Sheet "本地生活" (Local life) ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2838

# ---- Sheet "全部类型" (All types) ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 8757
$ws.Range("F8").Value = 16
$ws.Range("F9").Value = 11165
$ws.Range("F19").Value = 1903
$ws.Range("F20").Value = 710
$ws.Range("F21").Value = 635
$ws.Range("F22").Value = 359
$ws.Range("F29").Value = 1303
$ws.Range("F32").Value = 13
$ws.Range("F33").Value = 14
$ws.Range("F41").Value = 389
$ws.Range("F42").Value = 118
$ws.Range("F45").Value = 395
$ws.Range("F48").Value = 158
$ws.Range("F49").Value = 144
